$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 205.90909
$ws.Range("I6").Value = 205.90909
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 617.72727
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -505.72727
$ws.Range("N6").ClearContents()

$ws.Range("H40").Value = 3977.9
$ws.Range("I40").Value = 2716
$ws.Range("K40").Value = 2716
$ws.Range("M40").Value = -2541

$ws.Range("H48").Value = 3397.8
$ws.Range("J48").Value = 3397.8
$ws.Range("L48").Value = 10193.4
$ws.Range("N48").Value = -10777.4

$ws.Range("H56").Value = 3397.8
$ws.Range("J56").Value = 3397.8
$ws.Range("L56").Value = 10193.4
$ws.Range("N56").Value = -11261.4

$ws.Range("H63").Value = 69999
$ws.Range("J63").Value = 69999
$ws.Range("L63").Value = 69999
$ws.Range("N63").Value = -71247

$ws.Range("H66").Value = 69999
$ws.Range("J66").Value = 69999
$ws.Range("L66").Value = 209997
$ws.Range("N66").Value = -216237

$ws.Range("H98").Value = 1410.1562
$ws.Range("I98").Value = 1181.4517
$ws.Range("K98").Value = 1181.4517
$ws.Range("M98").Value = 316.5482999999999

$ws.Range("H113").Value = 5998.25
$ws.Range("I113").Value = 5998
$ws.Range("K113").Value = 5998
$ws.Range("M113").Value = -2744

$ws.Range("H116").Value = 6137
$ws.Range("J116").Value = 6232.875
$ws.Range("L116").Value = 6232.875
$ws.Range("N116").Value = -13116.875

$ws.Range("H122").Value = 1410.1562
$ws.Range("I122").Value = 1181.4517
$ws.Range("K122").Value = 3544.3551
$ws.Range("M122").Value = -1094.3551

$ws.Range("H138").Value = 2812.9768
$ws.Range("I138").Value = 2701.3667
$ws.Range("J138").Value = 3070.5386
$ws.Range("K138").Value = 8104.1001
$ws.Range("L138").Value = 9211.6158
$ws.Range("M138").Value = -2964.1001
$ws.Range("N138").Value = -19491.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3038.7273
$ws.Range("I2").Value = 3038.7273
$ws.Range("K2").Value = 3038.7273
$ws.Range("M2").Value = -2925.7273

$ws.Range("H14").Value = 7745
$ws.Range("I14").Value = 666
$ws.Range("J14").Value = 10104.667
$ws.Range("K14").Value = 666
$ws.Range("L14").Value = 10104.667
$ws.Range("M14").Value = -491
$ws.Range("N14").Value = -10454.667

$ws.Range("H43").Value = 34823.5
$ws.Range("I43").Value = 34648
$ws.Range("J43").Value = 34999
$ws.Range("K43").Value = 34648
$ws.Range("L43").Value = 34999
$ws.Range("M43").Value = -34335
$ws.Range("N43").Value = -35625

$ws.Range("H61").Value = 3496.48
$ws.Range("I61").Value = 3475.5
$ws.Range("K61").Value = 3475.5
$ws.Range("M61").Value = -3263.5

$ws.Range("H108").Value = 25000
$ws.Range("J108").Value = 25000
$ws.Range("L108").Value = 25000
$ws.Range("N108").Value = -32680

$ws.Range("H116").Value = 3038.7273
$ws.Range("I116").Value = 3038.7273
$ws.Range("K116").Value = 3038.7273
$ws.Range("M116").Value = -744.7273

$ws.Range("H122").Value = 6099.375
$ws.Range("I122").Value = 5399.3335
$ws.Range("K122").Value = 16198.0005
$ws.Range("M122").Value = -13748.0005

$ws.Range("H136").Value = 3496.48
$ws.Range("I136").Value = 3475.5
$ws.Range("K136").Value = 10426.5
$ws.Range("M136").Value = -7876.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3038.7273
$ws.Range("I3").Value = 3038.7273
$ws.Range("K3").Value = 3038.7273
$ws.Range("M3").Value = -2924.7273

$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()

$ws.Range("H115").Value = 84998
$ws.Range("J115").Value = 84998
$ws.Range("L115").Value = 84998
$ws.Range("N115").Value = -88132

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1393.5385
$ws.Range("I16").Value = 1393.5385
$ws.Range("K16").Value = 1393.5385
$ws.Range("M16").Value = -1106.5385

$ws.Range("H31").Value = 6851.706
$ws.Range("I31").Value = 5999.1113
$ws.Range("J31").Value = 7810.875
$ws.Range("K31").Value = 5999.1113
$ws.Range("L31").Value = 7810.875
$ws.Range("M31").Value = -5704.1113
$ws.Range("N31").Value = -8400.875

$ws.Range("H34").Value = 6851.706
$ws.Range("I34").Value = 5999.1113
$ws.Range("J34").Value = 7810.875
$ws.Range("K34").Value = 5999.1113
$ws.Range("L34").Value = 7810.875
$ws.Range("M34").Value = -5797.1113
$ws.Range("N34").Value = -8214.875

$ws.Range("H58").Value = 40353.04
$ws.Range("I58").Value = 57099.445
$ws.Range("K58").Value = 57099.445
$ws.Range("M58").Value = -56896.445

$ws.Range("H99").Value = 4281.9
$ws.Range("I99").Value = 4443.8
$ws.Range("K99").Value = 4443.8
$ws.Range("M99").Value = -2945.8

$ws.Range("H113").Value = 1393.5385
$ws.Range("I113").Value = 1393.5385
$ws.Range("K113").Value = 1393.5385
$ws.Range("M113").Value = 776.4614999999999

$ws.Range("H122").Value = 2699.8572
$ws.Range("I122").Value = 2699.8572
$ws.Range("K122").Value = 8099.571599999999
$ws.Range("M122").Value = -5649.571599999999

$ws.Range("H126").Value = 4281.9
$ws.Range("I126").Value = 4443.8
$ws.Range("K126").Value = 13331.4
$ws.Range("M126").Value = -10861.4

$ws.Range("H134").Value = 49083.41
$ws.Range("I134").Value = 61196.707
$ws.Range("J134").Value = 7898.2
$ws.Range("K134").Value = 183590.121
$ws.Range("L134").Value = 23694.6
$ws.Range("M134").Value = -181055.121
$ws.Range("N134").Value = -28764.6

$ws.Range("H136").Value = 40353.04
$ws.Range("I136").Value = 57099.445
$ws.Range("K136").Value = 171298.335
$ws.Range("M136").Value = -168748.335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 271626.6
$ws.Range("I4").Value = 290992.78
$ws.Range("K4").Value = 872978.3400000001
$ws.Range("M4").Value = -872866.3400000001

$ws.Range("H11").Value = 657.7143
$ws.Range("I11").Value = 99.5
$ws.Range("J11").Value = 881
$ws.Range("K11").Value = 298.5
$ws.Range("L11").Value = 2643
$ws.Range("M11").Value = -158.5
$ws.Range("N11").Value = -2923

$ws.Range("H13").Value = 33350.332
$ws.Range("I13").Value = 33350.332
$ws.Range("K13").Value = 100050.996
$ws.Range("M13").Value = -99882.99600000001

$ws.Range("H105").Value = 9082.25
$ws.Range("J105").Value = 9082.25
$ws.Range("L105").Value = 27246.75
$ws.Range("N105").Value = -32488.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 23638.334
$ws.Range("J92").Value = 23638.334
$ws.Range("L92").Value = 23638.334
$ws.Range("N92").Value = -27382.334

$ws.Range("H102").Value = 2895.2917
$ws.Range("I102").Value = 2294.2104
$ws.Range("K102").Value = 2294.2104
$ws.Range("M102").Value = -672.2103999999999

$ws.Range("H108").Value = 67499.5
$ws.Range("J108").Value = 67499.5
$ws.Range("L108").Value = 67499.5
$ws.Range("N108").Value = -75179.5

$ws.Range("H132").Value = 162573.05
$ws.Range("I132").Value = 162573.05
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 487719.15
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -485189.15
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5209.2666
$ws.Range("I40").Value = 4413.2
$ws.Range("K40").Value = 4413.2
$ws.Range("M40").Value = -4277.2

$ws.Range("H61").Value = 5854.067
$ws.Range("I61").Value = 4562.3076
$ws.Range("K61").Value = 4562.3076
$ws.Range("M61").Value = -4360.3076

$ws.Range("H113").Value = 5854.067
$ws.Range("I113").Value = 4562.3076
$ws.Range("K113").Value = 4562.3076
$ws.Range("M113").Value = -2392.3076

$ws.Range("H122").Value = 4926.278
$ws.Range("I122").Value = 4397.875
$ws.Range("K122").Value = 13193.625
$ws.Range("M122").Value = -10743.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6604.0557
$ws.Range("I122").Value = 8124.923
$ws.Range("J122").Value = 2649.8
$ws.Range("K122").Value = 24374.769
$ws.Range("L122").Value = 7949.400000000001
$ws.Range("M122").Value = -21924.769
$ws.Range("N122").Value = -12849.4

$ws.Range("H132").Value = 79856.84
$ws.Range("I132").Value = 86362
$ws.Range("J132").Value = 1795
$ws.Range("K132").Value = 259086
$ws.Range("L132").Value = 5385
$ws.Range("M132").Value = -256556
$ws.Range("N132").Value = -10445
